# Updated cryptos list on Mon May 15 13:46:36 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for every
# coin row, and reflects the EnergySwap/Quant rows trading ranking places
# (row 45 <-> row 46), including their Coin name + Link columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds values such as "27.717.23" or "0.4314" that must
# stay plain text (Excel would otherwise silently reinterpret anything that
# parses as a number, e.g. "1.013" -> 1.013 numeric, dropping the trailing
# zero / changing the stored type). Force the whole Price/Volume block to
# Text format first, write the literal strings, then clear the formatting
# back off again so the cells end up with no explicit style - exactly like
# the untouched cells around them.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

function Set-Row([int]$row, [string]$price, [string]$volume) {
    $ws.Cells.Item($row, 4).Value = $price
    $ws.Cells.Item($row, 5).Value = $volume
}

Set-Row 2  "27.740.10"    "  -0.05%  "
Set-Row 3  "1.848.65"     "  -0.83%  "
$ws.Cells.Item(4, 5).Value  = "  -2.03%  "
Set-Row 5  "320.29"       "  -1.25%  "
$ws.Cells.Item(6, 5).Value  = "  -2.06%  "
Set-Row 7  "0.4315"       "  -2.45%  "
Set-Row 8  "0.3743"       "  -1.48%  "
Set-Row 9  "0.07368"      "  -1.22%  "
Set-Row 10 "0.8790"       "  -0.74%  "
Set-Row 11 "21.70"        "  -0.27%  "
Set-Row 12 "1.858.93"     "  -0.45%  "
Set-Row 13 "6.738"        "  -0.44%  "
Set-Row 14 "5.459"        "  -1.69%  "
Set-Row 15 "0.07135"      "  -1.10%  "
Set-Row 16 "88.27"        "  +4.63%  "
$ws.Cells.Item(17, 5).Value = "  -2.17%  "
Set-Row 18 "0.000009008"  "  -1.13%  "
$ws.Cells.Item(19, 5).Value = "  -2.03%  "
Set-Row 20 "15.48"        "  -0.61%  "
Set-Row 21 "27.748.01"    "  -0.03%  "
Set-Row 22 "5.229"        "  -1.49%  "
Set-Row 23 "11.12"        "  -1.96%  "
Set-Row 24 "2.085.49"     "  -0.34%  "
Set-Row 25 "2.015"        "  -0.14%  "
Set-Row 26 "155.72"       "  -1.93%  "
$ws.Cells.Item(27, 5).Value = "  -1.10%  "
Set-Row 28 "2.148"        "  +7.66%  "
$ws.Cells.Item(29, 5).Value = "  +1.27%  "
Set-Row 30 "119.11"       "  +0.77%  "
Set-Row 31 "0.08963"      "  -0.81%  "
Set-Row 32 "1.234"        "  +0.98%  "
Set-Row 33 "0.7806"       "  +0.21%  "
$ws.Cells.Item(34, 5).Value = "  -0.17%  "
$ws.Cells.Item(35, 5).Value = "  -3.49%  "
Set-Row 36 "1.013"        "  -2.17%  "
$ws.Cells.Item(37, 5).Value = "  -1.23%  "
$ws.Cells.Item(38, 5).Value = "  +0.13%  "
$ws.Cells.Item(39, 5).Value = "  -0.94%  "
Set-Row 40 "7.307"        "  +6.28%  "
Set-Row 41 "2.881"        "  +0.14%  "
Set-Row 42 "0.5146"       "  -1.15%  "
Set-Row 43 "0.1690"       "  -0.24%  "
Set-Row 44 "8.826"        "  +1.68%  "

# Rows 45/46 swapped ranking places: EnergySwap moved from #45 to #46 and
# Quant moved from #46 to #45; the Coin name and Link columns move with them
# while Price/Volume reflect the latest market data for each coin.
$ws.Cells.Item(45, 2).Value = "Quant"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-Row 45 "109.71" "  -0.75%  "

$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-Row 46 "10.72" "  +0.27%  "

Set-Row 47 "0.4756"  "  +0.82%  "
Set-Row 48 "0.06494" "  -3.73%  "
$ws.Cells.Item(49, 4).Value = "1.699"
Set-Row 50 "1.013"   "  -2.27%  "
Set-Row 51 "1.860"   "  -2.90%  "

$dataRange.ClearFormats()
